# daily auto push: 2026-01-09 09:39 UTC
#
# The daily scraper appends a new "time seen" reading for the most recent
# date/weekday group. The new reading (2026/01/09, 金, time 17) belongs
# right after the existing 2026/01/09 rows (596-598), so every row from
# the old row 599 onward needs to shift down by one to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 599 - this pushes the old rows 599..640 down
# to 600..641 and widens the sheet's used range automatically.
$ws.Rows.Item(599).Insert()

# Column A holds the date as plain text (e.g. "2026/01/09"), not a real
# date serial. A bare string assignment gets auto-converted to a date by
# the Value setter, so prefix it with an apostrophe to force text entry,
# then reset the cell style back to Normal so no stray date number
# format is left behind on the cell.
$ws.Range("A599").Value = "'2026/01/09"
$ws.Range("A599").Style = "Normal"

$ws.Range("B599").Value = "金"
$ws.Range("C599").Value = 17
$ws.Range("D599").Value = 201
